$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Color" and "Semester" columns (D and E) are swapped: Semester now
# comes right after ECTS (column D), and Color moves to column E.

# Header row
$ws.Range("D1").Value = "Semester"
$ws.Range("E1").Value = "Color"

# Data rows: D gets the (numeric) semester value, E gets the (text) color value
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = "RED"

$ws.Range("D3").Value = 1
$ws.Range("E3").Value = "RED"

$ws.Range("D4").Value = 1
$ws.Range("E4").Value = "RED"

$ws.Range("D5").Value = 2
$ws.Range("E5").Value = "BLUE"

$ws.Range("D6").Value = 3
$ws.Range("E6").Value = "BLUE"

# Selection moved from E7 to G4
$ws.Range("G4").Select()
